$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 3) that mirrors row 2's "ficha" entry, with DIA = 10
$ws.Range("A3").Value = "CC"
$ws.Range("B3").Value = 1070593778
$ws.Range("C3").Value = "MARLON"
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = "MAYO"
$ws.Range("F3").Value = 2024
$ws.Range("G3").Value = 2671143

$ws.Range("G3").Select()
